# Append 45 new master-data rows (102-146) to the
# "master-reg_center_machine_devic" sheet, matching the pattern of the
# existing rows (A/B cycle through a 9-value block, C increments by 1,
# D/E/F/G are constant: "eng" / TRUE / "superadmin()" / "now()").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=regcntr_id, B=machine_id, C=device_id
$newRows = @(
    @(10002, 10021, 3000121),
    @(10003, 10022, 3000122),
    @(10004, 10023, 3000123),
    @(10005, 10024, 3000124),
    @(10006, 10025, 3000125),
    @(10007, 10026, 3000126),
    @(10008, 10027, 3000127),
    @(10009, 10028, 3000128),
    @(10010, 10029, 3000129),
    @(10002, 10021, 3000130),
    @(10003, 10022, 3000131),
    @(10004, 10023, 3000132),
    @(10005, 10024, 3000133),
    @(10006, 10025, 3000134),
    @(10007, 10026, 3000135),
    @(10008, 10027, 3000136),
    @(10009, 10028, 3000137),
    @(10010, 10029, 3000138),
    @(10002, 10021, 3000139),
    @(10003, 10022, 3000140),
    @(10004, 10023, 3000141),
    @(10005, 10024, 3000142),
    @(10006, 10025, 3000143),
    @(10007, 10026, 3000144),
    @(10008, 10027, 3000145),
    @(10009, 10028, 3000146),
    @(10010, 10029, 3000147),
    @(10002, 10021, 3000148),
    @(10003, 10022, 3000149),
    @(10004, 10023, 3000150),
    @(10005, 10024, 3000151),
    @(10006, 10025, 3000152),
    @(10007, 10026, 3000153),
    @(10008, 10027, 3000154),
    @(10009, 10028, 3000155),
    @(10010, 10029, 3000156),
    @(10002, 10021, 3000157),
    @(10003, 10022, 3000158),
    @(10004, 10023, 3000159),
    @(10005, 10024, 3000160),
    @(10006, 10025, 3000161),
    @(10007, 10026, 3000162),
    @(10008, 10027, 3000163),
    @(10009, 10028, 3000164),
    @(10010, 10029, 3000165)
)

$startRow = 102
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the author's final selection/view state: the newly added block
# A102:G146 is selected with A102 as the active cell.
[void]$ws.Range("A102:G146").Select()

# The saved workbook also carries an explicit portrait page setup.
$ws.PageSetup.Orientation = 1

Write-Output "Added $($newRows.Count) rows (102-146) to $($ws.Name())"
